$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.134.43'
$ws.Range("E2").Value = '  -1.96%  '

$ws.Range("D3").Value = '2.292.39'
$ws.Range("E3").Value = '  -3.40%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.62'
$ws.Range("E5").Value = '  -0.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.02'
$ws.Range("E6").Value = '  -4.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  -1.18%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  -3.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.33'
$ws.Range("E10").Value = '  -6.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  -2.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.24'
$ws.Range("E12").Value = '  -3.89%  '

$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("E14").Value = '  -5.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.23'
$ws.Range("E15").Value = '  -6.04%  '

$ws.Range("D16").Value = '2.638.29'

$ws.Range("D17").Value = '2.296.60'
$ws.Range("E17").Value = '  -3.51%  '

$ws.Range("D18").Value = '42.009.21'
$ws.Range("E18").Value = '  -2.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("E19").Value = '  -2.52%  '

$ws.Range("E20").Value = '  -1.84%  '

$ws.Range("B21").Value = 'PancakeSwap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.65'
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.22'
$ws.Range("E22").Value = '  -4.13%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '281.77'
$ws.Range("E23").Value = '  +9.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.02'
$ws.Range("E24").Value = '  +5.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  -3.54%  '

$ws.Range("E26").Value = '  +0.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.79'
$ws.Range("E27").Value = '  -6.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.34'
$ws.Range("E28").Value = '  +3.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.92'
$ws.Range("E29").Value = '  -0.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.56'
$ws.Range("E30").Value = '  -1.42%  '

$ws.Range("E31").Value = '  -4.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0873'
$ws.Range("E32").Value = '  -3.40%  '

$ws.Range("E33").Value = '  -3.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.82'
$ws.Range("E34").Value = '  -3.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.137'
$ws.Range("E35").Value = '  +4.00%  '

$ws.Range("E36").Value = '  -6.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  -4.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.91'
$ws.Range("E38").Value = '  +7.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0349'
$ws.Range("E39").Value = '  -4.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.71'
$ws.Range("E40").Value = '  -4.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.50'
$ws.Range("E41").Value = '  +11.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.46'
$ws.Range("E42").Value = '  -6.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.35'
$ws.Range("E43").Value = '  -3.08%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("E45").Value = '  -7.45%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.99'
$ws.Range("E46").Value = '  -3.35%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.30'
$ws.Range("E47").Value = '  +0.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '77.14'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.96'
$ws.Range("E49").Value = '  -3.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.27'
$ws.Range("E50").Value = '  -5.36%  '

$ws.Range("D51").Value = '1.564.85'
$ws.Range("E51").Value = '  -0.60%  '
